$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to touch to Text format so that
# numeric-looking strings (e.g. "588.92") are stored as text, matching the
# inlineStr cell type used throughout column D - then restore the default
# "Normal" style so no stray style index/attribute is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "64.135.95"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "2.782.70"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "588.92"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "160.58"
$ws.Range("E6").Value = "  +7.29%  "
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  +2.05%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").Value = "0.398"
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").Value = "3.278.95"
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("D14").Value = "27.44"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").Value = "64.065.19"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "0.0000159"
$ws.Range("E16").Value = "  +5.68%  "
$ws.Range("D17").Value = "2.789.10"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "12.41"
$ws.Range("E18").Value = "  +3.95%  "
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("D20").Value = "366.92"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "7.06"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "0.574"
$ws.Range("E22").Value = "  +7.35%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "67.55"
$ws.Range("E24").Value = "  +3.12%  "
$ws.Range("E25").Value = "  +6.11%  "
$ws.Range("D26").Value = "8.79"
$ws.Range("E26").Value = "  +3.09%  "
$ws.Range("D27").Value = "0.0₃0969"
$ws.Range("E27").Value = "  +12.77%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("D31").Value = "1.28"
$ws.Range("E31").Value = "  +7.98%  "
$ws.Range("E32").Value = "  +10.27%  "
$ws.Range("D33").Value = "172.04"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").Value = "20.88"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E36").Value = "  +5.87%  "
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").Value = "4.28"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").Value = "341.82"
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("E41").Value = "  +11.56%  "
$ws.Range("D42").Value = "40.30"
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("D43").Value = "22.52"
$ws.Range("E43").Value = "  +4.14%  "
$ws.Range("D44").Value = "22.51"
$ws.Range("E44").Value = "  +4.00%  "
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").Value = "0.0262"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("D48").Value = "138.78"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("D50").Value = "2.171.76"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("E51").Value = "  +0.39%  "

$priceRange.Style = "Normal"
